# Swap the betting-odds data between specific pairs of rows, while keeping
# the "id" column (A) and "Date" column (D) anchored to their original row
# position. Columns B, C and E..AD are swapped between each pair of rows.
#
# Row pairs that need to be swapped (1-based worksheet row numbers):
$pairs = @(
    @(47, 48),
    @(71, 72),
    @(101, 102),
    @(109, 110),
    @(149, 150),
    @(229, 231),
    @(232, 233),
    @(263, 265),
    @(310, 311)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    # Swap columns B:C
    $rangeA_BC = $ws.Range("B$rowA", "C$rowA")
    $rangeB_BC = $ws.Range("B$rowB", "C$rowB")
    $valA_BC = $rangeA_BC.Value2
    $valB_BC = $rangeB_BC.Value2
    $rangeA_BC.Value2 = $valB_BC
    $rangeB_BC.Value2 = $valA_BC

    # Swap columns E:AD (skip A = id and D = Date, which stay with the row)
    $rangeA_EAD = $ws.Range("E$rowA", "AD$rowA")
    $rangeB_EAD = $ws.Range("E$rowB", "AD$rowB")
    $valA_EAD = $rangeA_EAD.Value2
    $valB_EAD = $rangeB_EAD.Value2
    $rangeA_EAD.Value2 = $valB_EAD
    $rangeB_EAD.Value2 = $valA_EAD
}
